$wb = $excel.ActiveWorkbook

function Set-CellValue($ws, $ref, $val) {
    $ws.Range($ref).Value = $val
}

function Clear-CellValue($ws, $ref) {
    $ws.Range($ref).ClearContents()
}

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws "H17" 535.3953
Set-CellValue $ws "J17" 535.3953
Set-CellValue $ws "L17" 1606.1859
Set-CellValue $ws "N17" -1942.1859
Set-CellValue $ws "H33" 410.82858
Set-CellValue $ws "I33" 147.43478
Set-CellValue $ws "J33" 915.6667
Set-CellValue $ws "K33" 147.43478
Set-CellValue $ws "L33" 915.6667
Set-CellValue $ws "M33" 81.56522000000001
Set-CellValue $ws "N33" -1373.6667
Set-CellValue $ws "H70" 3786.3635
Set-CellValue $ws "J70" 3530
Set-CellValue $ws "L70" 10590
Set-CellValue $ws "N70" -11130
Set-CellValue $ws "H73" 3786.3635
Set-CellValue $ws "J73" 3530
Set-CellValue $ws "L73" 10590
Set-CellValue $ws "N73" -12462
Set-CellValue $ws "H113" 3980
Set-CellValue $ws "I113" 3828.5715
Set-CellValue $ws "J113" 4333.3335
Set-CellValue $ws "K113" 3828.5715
Set-CellValue $ws "L113" 4333.3335
Set-CellValue $ws "M113" -574.5715
Set-CellValue $ws "N113" -10841.3335
Set-CellValue $ws "H127" 1182.8182
Set-CellValue $ws "I127" 632.3333
Set-CellValue $ws "J127" 1843.4
Set-CellValue $ws "K127" 1896.9999
Set-CellValue $ws "L127" 5530.200000000001
Set-CellValue $ws "M127" 3063.0001
Set-CellValue $ws "N127" -15450.2
Set-CellValue $ws "H137" 2566374.8
Set-CellValue $ws "I137" 4001653.8
Set-CellValue $ws "J137" 3376.6428
Set-CellValue $ws "K137" 12004961.4
Set-CellValue $ws "L137" 10129.9284
Set-CellValue $ws "M137" -12002411.4
Set-CellValue $ws "N137" -15229.9284
Set-CellValue $ws "H138" 3273603
Set-CellValue $ws "I138" 287255.28
Set-CellValue $ws "J138" 8774770
Set-CellValue $ws "K138" 861765.8400000001
Set-CellValue $ws "L138" 26324310
Set-CellValue $ws "M138" -856625.8400000001
Set-CellValue $ws "N138" -26334590

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws "H2" 3109.1
Set-CellValue $ws "I2" 3368.5
Set-CellValue $ws "J2" 2720
Set-CellValue $ws "K2" 3368.5
Set-CellValue $ws "L2" 2720
Set-CellValue $ws "M2" -3255.5
Set-CellValue $ws "N2" -2946
Set-CellValue $ws "H32" 1268.36
Set-CellValue $ws "I32" 969.7027
Set-CellValue $ws "J32" 2118.3845
Set-CellValue $ws "K32" 969.7027
Set-CellValue $ws "L32" 2118.3845
Set-CellValue $ws "M32" -682.7027
Set-CellValue $ws "N32" -2692.3845
Set-CellValue $ws "H61" 62625930
Set-CellValue $ws "I61" 100100536
Set-CellValue $ws "J61" 168252.33
Set-CellValue $ws "K61" 100100536
Set-CellValue $ws "L61" 168252.33
Set-CellValue $ws "M61" -100100324
Set-CellValue $ws "N61" -168676.33
Set-CellValue $ws "H97" 1839195.5
Set-CellValue $ws "I97" 2404857.8
Set-CellValue $ws "J97" 793.75
Set-CellValue $ws "K97" 2404857.8
Set-CellValue $ws "L97" 793.75
Set-CellValue $ws "M97" -2404361.8
Set-CellValue $ws "N97" -1785.75
Set-CellValue $ws "H116" 3109.1
Set-CellValue $ws "I116" 3368.5
Set-CellValue $ws "J116" 2720
Set-CellValue $ws "K116" 3368.5
Set-CellValue $ws "L116" 2720
Set-CellValue $ws "M116" -1074.5
Set-CellValue $ws "N116" -7308
Set-CellValue $ws "H132" 96299.5
Set-CellValue $ws "I132" 63473.562
Set-CellValue $ws "J132" 183835.33
Set-CellValue $ws "K132" 190420.686
Set-CellValue $ws "L132" 551505.99
Set-CellValue $ws "M132" -187890.686
Set-CellValue $ws "N132" -556565.99
Set-CellValue $ws "H136" 62625930
Set-CellValue $ws "I136" 100100536
Set-CellValue $ws "J136" 168252.33
Set-CellValue $ws "K136" 300301608
Set-CellValue $ws "L136" 504756.99
Set-CellValue $ws "M136" -300299058
Set-CellValue $ws "N136" -509856.99

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
Set-CellValue $ws "H3" 3109.1
Set-CellValue $ws "I3" 3368.5
Set-CellValue $ws "J3" 2720
Set-CellValue $ws "K3" 3368.5
Set-CellValue $ws "L3" 2720
Set-CellValue $ws "M3" -3254.5
Set-CellValue $ws "N3" -2948
Set-CellValue $ws "H94" 524.2353000000001
Set-CellValue $ws "I94" 475.46667
Set-CellValue $ws "K94" 475.46667
Set-CellValue $ws "M94" -24.46667000000002
Set-CellValue $ws "H99" 1327.88
Set-CellValue $ws "I99" 1424.4117
Set-CellValue $ws "J99" 1122.75
Set-CellValue $ws "K99" 1424.4117
Set-CellValue $ws "L99" 1122.75
Set-CellValue $ws "M99" 73.58829999999989
Set-CellValue $ws "N99" -4118.75
Set-CellValue $ws "H124" 0
Set-CellValue $ws "J124" 0
Set-CellValue $ws "L124" 0
Clear-CellValue $ws "N124"

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws "H31" 3625.88
Set-CellValue $ws "I31" 1856.6842
Set-CellValue $ws "J31" 9228.333000000001
Set-CellValue $ws "K31" 1856.6842
Set-CellValue $ws "L31" 9228.333000000001
Set-CellValue $ws "M31" -1561.6842
Set-CellValue $ws "N31" -9818.333000000001
Set-CellValue $ws "H34" 3625.88
Set-CellValue $ws "I34" 1856.6842
Set-CellValue $ws "J34" 9228.333000000001
Set-CellValue $ws "K34" 1856.6842
Set-CellValue $ws "L34" 9228.333000000001
Set-CellValue $ws "M34" -1654.6842
Set-CellValue $ws "N34" -9632.333000000001
Set-CellValue $ws "H64" 29644.2
Set-CellValue $ws "J64" 29644.2
Set-CellValue $ws "L64" 29644.2
Set-CellValue $ws "N64" -30140.2
Set-CellValue $ws "H67" 29644.2
Set-CellValue $ws "J67" 29644.2
Set-CellValue $ws "L67" 29644.2
Set-CellValue $ws "N67" -31360.2
Set-CellValue $ws "H94" 5413.5454
Set-CellValue $ws "I94" 15650
Set-CellValue $ws "J94" 1574.875
Set-CellValue $ws "K94" 15650
Set-CellValue $ws "L94" 1574.875
Set-CellValue $ws "M94" -15199
Set-CellValue $ws "N94" -2476.875
Set-CellValue $ws "H99" 4971.1763
Set-CellValue $ws "I99" 6882.6665
Set-CellValue $ws "J99" 3928.5454
Set-CellValue $ws "K99" 6882.6665
Set-CellValue $ws "L99" 3928.5454
Set-CellValue $ws "M99" -5384.6665
Set-CellValue $ws "N99" -6924.5454
Set-CellValue $ws "H122" 2355.75
Set-CellValue $ws "I122" 1982.909
Set-CellValue $ws "J122" 2671.2307
Set-CellValue $ws "K122" 5948.727000000001
Set-CellValue $ws "L122" 8013.6921
Set-CellValue $ws "M122" -3498.727000000001
Set-CellValue $ws "N122" -12913.6921
Set-CellValue $ws "H126" 4971.1763
Set-CellValue $ws "I126" 6882.6665
Set-CellValue $ws "J126" 3928.5454
Set-CellValue $ws "K126" 20647.9995
Set-CellValue $ws "L126" 11785.6362
Set-CellValue $ws "M126" -18177.9995
Set-CellValue $ws "N126" -16725.6362
Set-CellValue $ws "H134" 24189.914
Set-CellValue $ws "I134" 1316.6389
Set-CellValue $ws "J134" 99047.91
Set-CellValue $ws "K134" 3949.9167
Set-CellValue $ws "L134" 297143.73
Set-CellValue $ws "M134" -1414.9167
Set-CellValue $ws "N134" -302213.73
Set-CellValue $ws "H141" 50548.223
Set-CellValue $ws "J141" 50548.223
Set-CellValue $ws "L141" 50548.223
Set-CellValue $ws "N141" -60908.223

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
Set-CellValue $ws "H107" 578.89746
Set-CellValue $ws "I107" 1068.6666
Set-CellValue $ws "J107" 361.22223
Set-CellValue $ws "K107" 3205.9998
Set-CellValue $ws "L107" 1083.66669
Set-CellValue $ws "M107" -1285.9998
Set-CellValue $ws "N107" -4923.66669
Set-CellValue $ws "H120" 13833.167
Set-CellValue $ws "I120" 7999
Set-CellValue $ws "J120" 15000
Set-CellValue $ws "K120" 23997
Set-CellValue $ws "L120" 45000
Set-CellValue $ws "M120" -19159
Set-CellValue $ws "N120" -54676
Set-CellValue $ws "H131" 10870572
Set-CellValue $ws "I131" 76923460
Set-CellValue $ws "J131" 1109.2405
Set-CellValue $ws "K131" 230770380
Set-CellValue $ws "L131" 3327.721500000001
Set-CellValue $ws "M131" -230765340
Set-CellValue $ws "N131" -13407.7215

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
Set-CellValue $ws "H70" 25410.396
Set-CellValue $ws "I70" 35546.875
Set-CellValue $ws "J70" 5137.4375
Set-CellValue $ws "K70" 35546.875
Set-CellValue $ws "L70" 5137.4375
Set-CellValue $ws "M70" -35276.875
Set-CellValue $ws "N70" -5677.4375
Set-CellValue $ws "H73" 25410.396
Set-CellValue $ws "I73" 35546.875
Set-CellValue $ws "J73" 5137.4375
Set-CellValue $ws "K73" 35546.875
Set-CellValue $ws "L73" 5137.4375
Set-CellValue $ws "M73" -34610.875
Set-CellValue $ws "N73" -7009.4375
Set-CellValue $ws "H80" 3292.647
Set-CellValue $ws "I80" 2278
Set-CellValue $ws "J80" 3715.4167
Set-CellValue $ws "K80" 2278
Set-CellValue $ws "L80" 3715.4167
Set-CellValue $ws "M80" -1280
Set-CellValue $ws "N80" -5711.4167
Set-CellValue $ws "H83" 3292.647
Set-CellValue $ws "I83" 2278
Set-CellValue $ws "J83" 3715.4167
Set-CellValue $ws "K83" 11390
Set-CellValue $ws "L83" 18577.0835
Set-CellValue $ws "M83" -6398
Set-CellValue $ws "N83" -28561.0835
Set-CellValue $ws "H141" 32850
Set-CellValue $ws "J141" 32850
Set-CellValue $ws "L141" 32850
Set-CellValue $ws "N141" -43210

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
Set-CellValue $ws "H40" 2867.3333
Set-CellValue $ws "I40" 2726
Set-CellValue $ws "J40" 3150
Set-CellValue $ws "K40" 2726
Set-CellValue $ws "L40" 3150
Set-CellValue $ws "M40" -2590
Set-CellValue $ws "N40" -3422
Set-CellValue $ws "H100" 1368.2572
Set-CellValue $ws "I100" 1121.2609
Set-CellValue $ws "K100" 1121.2609
Set-CellValue $ws "M100" -580.2609

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
Set-CellValue $ws "H46" 47500
Set-CellValue $ws "I46" 0
Set-CellValue $ws "J46" 47500
Set-CellValue $ws "K46" 0
Set-CellValue $ws "L46" 47500
Clear-CellValue $ws "M46"
Set-CellValue $ws "N46" -47962
Set-CellValue $ws "H113" 1291.3684
Set-CellValue $ws "I113" 885.9167
Set-CellValue $ws "J113" 1986.4286
Set-CellValue $ws "K113" 2657.7501
Set-CellValue $ws "L113" 5959.2858
Set-CellValue $ws "M113" -487.7501000000002
Set-CellValue $ws "N113" -10299.2858
Set-CellValue $ws "H126" 1250.421
Set-CellValue $ws "I126" 1000
Set-CellValue $ws "J126" 1594.75
Set-CellValue $ws "K126" 3000
Set-CellValue $ws "L126" 4784.25
Set-CellValue $ws "M126" -530
Set-CellValue $ws "N126" -9724.25
Set-CellValue $ws "H134" 47500
Set-CellValue $ws "I134" 0
Set-CellValue $ws "J134" 47500
Set-CellValue $ws "K134" 0
Set-CellValue $ws "L134" 142500
Clear-CellValue $ws "M134"
Set-CellValue $ws "N134" -147570
Set-CellValue $ws "H139" 0
Set-CellValue $ws "J139" 0
Set-CellValue $ws "L139" 0
Clear-CellValue $ws "N139"
Set-CellValue $ws "H140" 38572.57
Set-CellValue $ws "J140" 38572.57
Set-CellValue $ws "L140" 38572.57
Set-CellValue $ws "N140" -48932.57
